$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updates to time module: refresh Ticket Sales (Q) and Embarking (R) values
$updates = @(
    @{ Row = 3;   Q = 57;  R = 34 },
    @{ Row = 10;  Q = 25;  R = 11 },
    @{ Row = 17;  Q = 43;  R = 11 },
    @{ Row = 23;  Q = 7;   R = 3  },
    @{ Row = 32;  Q = 46;  R = 34 },
    @{ Row = 40;  Q = 36;  R = 22 },
    @{ Row = 49;  Q = 25;  R = 15 },
    @{ Row = 58;  Q = 57;  R = 45 },
    @{ Row = 66;  Q = 98;  R = 57 },
    @{ Row = 74;  Q = 13;  R = 6  },
    @{ Row = 78;  Q = 41;  R = 16 },
    @{ Row = 89;  Q = 55;  R = 18 },
    @{ Row = 97;  Q = 75;  R = 58 },
    @{ Row = 106; Q = 50 },
    @{ Row = 115; Q = 93;  R = 51 },
    @{ Row = 124; Q = 61;  R = 10 },
    @{ Row = 133; Q = 100; R = 62 },
    @{ Row = 142; Q = 88;  R = 40 }
)

foreach ($u in $updates) {
    $ws.Range("Q$($u.Row)").Value = $u.Q
    if ($u.ContainsKey("R")) {
        $ws.Range("R$($u.Row)").Value = $u.R
    }
}
